{"js": "// Word JS API (Office.js) edit script.\n// Body of: async (context) => { ... }\n//\n// Change: the paragraph that explains the \"Justifica la relaci\u00f3n...\" activity\n// had its trailing clause reworded:\n//   old: \"...m\u00e1ximo 1 cuartilla cuyas caracter\u00edsticas deben incluir Arial 12,\n//         interlineado 1.5 y texto justificado y m\u00e1rgenes normales. \"\n//   new: \"...m\u00e1ximo 1 cuartilla. El texto debe cumplir con las siguientes\n//         caracter\u00edsticas: letra Arial 12, interlineado 1.5, texto\n//         justificado y m\u00e1rgenes normales. \"\n// The leading part of the sentence (\"Justifica la relaci\u00f3n ... m\u00e1ximo 1\n// cuartilla\") is unchanged; only the trailing clause is replaced. We locate\n// that trailing clause with a search (unique in the document) and replace it\n// in place so the run keeps its original character formatting\n// (Arial, sz 22, lang es-MX).\n\nconst oldTail =\n  \" cuyas caracter\u00edsticas deben incluir Arial 12, interlineado 1.5 y texto justificado y m\u00e1rgenes normales. \";\nconst newTail =\n  \". El texto debe cumplir con las siguientes caracter\u00edsticas: letra Arial 12, interlineado 1.5, texto justificado y m\u00e1rgenes normales. \";\n\nconst results = context.document.body.search(oldTail, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target sentence tail not found in document body.\");\n}\n\n// Replace the matched range's text in place; the new text inherits the\n// character formatting (font, size, language) already applied to the\n// matched range, same as the rest of the sentence.\nresults.items[0].insertText(newTail, Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# Word COM interop edit script.\n# $word.ActiveDocument is already open as $d below.\n#\n# Change: the paragraph that explains the \"Justifica la relaci\u00f3n...\" activity\n# had its trailing clause reworded:\n#   old: \"...m\u00e1ximo 1 cuartilla cuyas caracter\u00edsticas deben incluir Arial 12,\n#         interlineado 1.5 y texto justificado y m\u00e1rgenes normales. \"\n#   new: \"...m\u00e1ximo 1 cuartilla. El texto debe cumplir con las siguientes\n#         caracter\u00edsticas: letra Arial 12, interlineado 1.5, texto\n#         justificado y m\u00e1rgenes normales. \"\n# The leading part of the sentence (\"Justifica la relaci\u00f3n ... m\u00e1ximo 1\n# cuartilla\") stays the same; only the trailing clause is swapped out. We\n# walk the paragraphs, find the one containing that sentence, and rewrite its\n# Range.Text in place (instead of using Find/Replace) so the run keeps its\n# original rsid/formatting attributes, matching how Word preserves a run when\n# only its interior text changes.\n\n$d = $word.ActiveDocument\n\n$oldTail = \" cuyas caracter\u00edsticas deben incluir Arial 12, interlineado 1.5 y texto justificado y m\u00e1rgenes normales. \"\n$newTail = \". El texto debe cumplir con las siguientes caracter\u00edsticas: letra Arial 12, interlineado 1.5, texto justificado y m\u00e1rgenes normales. \"\n\nforeach ($p in $d.Paragraphs) {\n  $t = $p.Range.Text\n  if ($t -like \"*Justifica la relaci\u00f3n*$oldTail*\") {\n    $hasParaMark = $t.EndsWith(\"`r\")\n    if ($hasParaMark) {\n      $body = $t.Substring(0, $t.Length - 1)\n    } else {\n      $body = $t\n    }\n    $newBody = $body.Replace($oldTail, $newTail)\n    $p.Range.Text = $newBody\n  }\n}\n"}
